# "Part 1" (Erlang-B blocking calculator) -> "Part 2" (Erlang-C waiting-time calculator)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the numeric-looking literal to be stored as text (shared string),
    # matching the source workbook's convention of keeping all Inputs/Results
    # "Values" as text, then drop back to the default "Normal" style so no
    # stray per-cell number format sticks around.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 1 - headers
$ws.Range("A1").Value = "Inputs"
$ws.Range("B1").Value = "Values"

# Row 2 - Arrival Rate
$ws.Range("A2").Value = "Arrival Rate"
Set-TextValue $ws.Range("B2") "5.0"

# Row 3 - Service Rate
$ws.Range("A3").Value = "Service Rate"
Set-TextValue $ws.Range("B3") "6.0"

# Row 4 - P(W > 0) Less Than
$ws.Range("A4").Value = "P(W > 0) Less Than"
Set-TextValue $ws.Range("B4") "0.2"

# Row 5 - E(W) Less Than (new row)
$ws.Range("A5").Value = "E(W) Less Than"
Set-TextValue $ws.Range("B5") "0.5"

# Row 6 - Results header
$ws.Range("A6").Value = "Results"

# Row 7 - Number of Servers
$ws.Range("A7").Value = "Number of Servers"
Set-TextValue $ws.Range("B7") "3"

# Row 8 - E(S)
$ws.Range("A8").Value = "E(S)"
Set-TextValue $ws.Range("B8") "0.05555555555555556"

# Row 9 - E(N) (new row)
$ws.Range("A9").Value = "E(N)"
Set-TextValue $ws.Range("B9") "0.38461538461538464"

# Rename the sheet/tab: "Part 1" -> "Part 2"
$ws.Name = "Part 2"
